# Update "想去人数" (interest counter) values for two events that appear
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 298
$wsExhibit.Range("F4").Value = 1295

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 298
$wsAll.Range("F5").Value = 1295
